$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Insert a new row at position 14 ("Gacha3BrokenEnergys" / int / 3").
# All existing rows 14-19 (GoldBoxTurnMin, GoldBoxTurnMax, FirstGoldBox,
# MaxAnalysisLevel, MaxGuideQuestId, MaxBrokenEnergy) shift down to 15-20,
# and every formula reference to those rows is adjusted automatically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = "Gacha3BrokenEnergys"
$ws.Range("B14").Value = "int"
$ws.Range("D14").Value = 3

# ---------------------------------------------------------------------------
# Gacha3Events value: 10 -> 9  (row 11, column D)
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 9

# ---------------------------------------------------------------------------
# Rebuild the F2 "json blob" formula so it also folds in the newly added
# Gacha3BrokenEnergys row (A14) right before the FirstGoldBox row, which is
# now A17 after the insert above.
# ---------------------------------------------------------------------------
$formula = "=`"{`"`"`"&`nA2&`"`"`":`"&VLOOKUP(A2,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA8&`"`"`":`"&VLOOKUP(A8,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA9&`"`"`":`"&VLOOKUP(A9,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA10&`"`"`":`"&VLOOKUP(A10,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA11&`"`"`":`"&VLOOKUP(A11,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA12&`"`"`":`"&VLOOKUP(A12,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA13&`"`"`":`"&VLOOKUP(A13,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA14&`"`"`":`"&VLOOKUP(A14,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)&`",`"`"`"&`nA17&`"`"`":`"&VLOOKUP(A17,`$A:`$D,MATCH(`$D`$1,`$A`$1:`$D`$1,0),0)`n&`"}`""
$ws.Range("F2").Formula = $formula

# Setting a multi-line formula text makes the engine auto-fit row 2's height
# to the wrapped text; the source workbook keeps row 2 at the sheet default,
# so put it back the way it was.
$ws.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# Match the author's final selection (cell A15 was left selected).
# ---------------------------------------------------------------------------
[void]$ws.Range("A15").Select()
